$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header "name" -> "material_name"
$ws.Range("A1").Value = "material_name"

# Remove the "stock" column (column F) entirely, shifting nothing right of it
$ws.Columns("F").Delete()

# Update selection to match target workbook state
$ws.Range("F11").Select()
